# Update "想去人数" (F column) figures across the four worksheets to reflect
# the latest scrape output (gh-pages data refresh at commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1138
$ws1.Range("F5").Value = 48
$ws1.Range("F7").Value = 778
$ws1.Range("F14").Value = 935
$ws1.Range("F16").Value = 2046
$ws1.Range("F17").Value = 525
$ws1.Range("F18").Value = 8588
$ws1.Range("F19").Value = 794
$ws1.Range("F23").Value = 25
$ws1.Range("F26").Value = 1

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 527
$ws2.Range("F9").Value = 130
$ws2.Range("F12").Value = 15

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 423
$ws3.Range("F4").Value = 408

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 423
$ws4.Range("F5").Value = 408
$ws4.Range("F6").Value = 527
$ws4.Range("F7").Value = 1138
$ws4.Range("F10").Value = 48
$ws4.Range("F12").Value = 778
$ws4.Range("F23").Value = 935
$ws4.Range("F25").Value = 130
$ws4.Range("F27").Value = 2046
$ws4.Range("F28").Value = 525
$ws4.Range("F29").Value = 8589
$ws4.Range("F31").Value = 15
$ws4.Range("F32").Value = 794
$ws4.Range("F37").Value = 25
$ws4.Range("F45").Value = 1
